$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.020.29"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "'1.872.89"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.72%  "

$ws.Range("D5").Value = "'312.69"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("D7").Value = "'0.5038"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "'0.3830"
$ws.Range("E8").Value = "  -2.12%  "

$ws.Range("D9").Value = "'0.08979"
$ws.Range("E9").Value = "  -5.84%  "

$ws.Range("D10").Value = "'1.118"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("B11").Value = "Polkadot"
$ws.Range("C11").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D11").Value = "'6.378"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'20.70"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.862.72"
$ws.Range("E13").Value = "  -2.68%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.255"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'0.9995"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001105"
$ws.Range("E16").Value = "  -0.74%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'91.17"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06654"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'18.21"
$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.129"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "'28.030.37"
$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'11.50"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.264"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'3.395"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").Value = "'2.067.39"
$ws.Range("E26").Value = "  -3.44%  "

$ws.Range("D27").Value = "'2.500"
$ws.Range("E27").Value = "  -4.89%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.73"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'157.29"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "'126.67"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").Value = "'1.057"
$ws.Range("E32").Value = "  -1.98%  "

$ws.Range("D33").Value = "'5.609"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").Value = "'3.595"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").Value = "'9.470"
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("D36").Value = "'0.06597"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").Value = "'0.02404"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "'0.2191"
$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").Value = "'1.287"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("D40").Value = "'1.209"
$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("D41").Value = "'0.6385"
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("D42").Value = "'11.48"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").Value = "'4.905"
$ws.Range("E43").Value = "  -1.13%  "

$ws.Range("D44").Value = "'0.9990"
$ws.Range("E44").Value = "  -0.49%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.22"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6023"
$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("D47").Value = "'1.277"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("E48").Value = "  -1.95%  "

$ws.Range("D49").Value = "'1.240"
$ws.Range("E49").Value = "  +5.12%  "

$ws.Range("D50").Value = "'1.998"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").Value = "'120.81"
$ws.Range("E51").Value = "  -1.13%  "
